$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pretty")
$ws.Columns.Item(4).Insert()
Write-Host ("done")
